$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: translate Spanish column names to snake_case English field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the lowercase connector words ("de", "del", "la", "el") in a
# handful of state/municipality names so they match the rest of the data.
$ws.Range("A19").Value = "Ciudad De México"
$ws.Range("A25").Value = "Estado De México"
$ws.Range("B25").Value = "Ixtapan De La Sal"
$ws.Range("B26").Value = "San Felipe Del Progreso"
$ws.Range("B34").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B35").Value = "Silao De La Victoria"
$ws.Range("B38").Value = "Atoyac De Álvarez"
$ws.Range("B41").Value = "Tlapa De Comonfort"
$ws.Range("B44").Value = "Huejuquilla El Alto"
$ws.Range("B46").Value = "Unión De Tula"
$ws.Range("B59").Value = "Ocotlán De Morelos"
$ws.Range("B64").Value = "Jalpan De Serra"
$ws.Range("B67").Value = "Villa De La Paz"
$ws.Range("B82").Value = "Poza Rica De Hidalgo"

# Drop the trailing footnote/metadata rows (sample size, source, author,
# date) that used to live below the data table.
$ws.Rows("86:90").Delete()
